$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 (Severin Standgrill mit Grillplatte PG 8563), shifting all following rows up by one.
$ws.Rows.Item(3).Delete()

# Update the timestamp column (O) for the banner row and all product rows (2 through 30 after the shift)
# to reflect the new crawl time.
$ws.Range("O2:O30").Value = "2022-07-20 20:58:43"
